$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "48×23="
$t.Cell(1,2).Range.Text = "90×30="
$t.Cell(1,3).Range.Text = "97×76="
$t.Cell(1,4).Range.Text = "50×39="
$t.Cell(1,5).Range.Text = "19×29="
$t.Cell(2,1).Range.Text = "100×95="
$t.Cell(2,2).Range.Text = "52×41="
$t.Cell(2,3).Range.Text = "46×87="
$t.Cell(2,4).Range.Text = "82×84="
$t.Cell(2,5).Range.Text = "64×65="
$t.Cell(3,1).Range.Text = "95×30="
$t.Cell(3,2).Range.Text = "99×72="
$t.Cell(3,3).Range.Text = "58×24="
$t.Cell(3,4).Range.Text = "97×32="
$t.Cell(3,5).Range.Text = "80×68="
$t.Cell(4,1).Range.Text = "71×50="
$t.Cell(4,2).Range.Text = "76×40="
$t.Cell(4,3).Range.Text = "45×30="
$t.Cell(4,4).Range.Text = "39×68="
$t.Cell(4,5).Range.Text = "35×75="
$t.Cell(5,1).Range.Text = "52×80="
$t.Cell(5,2).Range.Text = "79×17="
$t.Cell(5,3).Range.Text = "84×21="
$t.Cell(5,4).Range.Text = "11×43="
$t.Cell(5,5).Range.Text = "72×54="
$t.Cell(6,1).Range.Text = "99×18="
$t.Cell(6,2).Range.Text = "17×84="
$t.Cell(6,3).Range.Text = "84×56="
$t.Cell(6,4).Range.Text = "12×16="
$t.Cell(6,5).Range.Text = "75×21="
$t.Cell(7,1).Range.Text = "47×87="
$t.Cell(7,2).Range.Text = "67×19="
$t.Cell(7,3).Range.Text = "48×62="
$t.Cell(7,4).Range.Text = "73×28="
$t.Cell(7,5).Range.Text = "98×86="
$t.Cell(8,1).Range.Text = "45×39="
$t.Cell(8,2).Range.Text = "79×72="
$t.Cell(8,3).Range.Text = "31×39="
$t.Cell(8,4).Range.Text = "68×25="
$t.Cell(8,5).Range.Text = "23×98="
$t.Cell(9,1).Range.Text = "62×72="
$t.Cell(9,2).Range.Text = "45×74="
$t.Cell(9,3).Range.Text = "46×12="
$t.Cell(9,4).Range.Text = "38×67="
$t.Cell(9,5).Range.Text = "70×16="
$t.Cell(10,1).Range.Text = "54×60="
$t.Cell(10,2).Range.Text = "29×27="
$t.Cell(10,3).Range.Text = "50×20="
$t.Cell(10,4).Range.Text = "47×14="
$t.Cell(10,5).Range.Text = "67×48="
$t.Cell(11,1).Range.Text = "13×66="
$t.Cell(11,2).Range.Text = "79×87="
$t.Cell(11,3).Range.Text = "80×22="
$t.Cell(11,4).Range.Text = "70×44="
$t.Cell(11,5).Range.Text = "50×28="
$t.Cell(12,1).Range.Text = "68×16="
$t.Cell(12,2).Range.Text = "93×14="
$t.Cell(12,3).Range.Text = "16×93="
$t.Cell(12,4).Range.Text = "78×20="
$t.Cell(12,5).Range.Text = "87×72="
$t.Cell(13,1).Range.Text = "58×48="
$t.Cell(13,2).Range.Text = "100×98="
$t.Cell(13,3).Range.Text = "40×97="
$t.Cell(13,4).Range.Text = "20×27="
$t.Cell(13,5).Range.Text = "29×20="
$t.Cell(14,1).Range.Text = "36×82="
$t.Cell(14,2).Range.Text = "75×50="
$t.Cell(14,3).Range.Text = "64×87="
$t.Cell(14,4).Range.Text = "39×30="
$t.Cell(14,5).Range.Text = "61×74="
$t.Cell(15,1).Range.Text = "66×69="
$t.Cell(15,2).Range.Text = "42×10="
$t.Cell(15,3).Range.Text = "69×56="
$t.Cell(15,4).Range.Text = "12×28="
$t.Cell(15,5).Range.Text = "41×49="
$t.Cell(16,1).Range.Text = "80×36="
$t.Cell(16,2).Range.Text = "65×76="
$t.Cell(16,3).Range.Text = "25×28="
$t.Cell(16,4).Range.Text = "90×69="
$t.Cell(16,5).Range.Text = "46×53="
$t.Cell(17,1).Range.Text = "62×12="
$t.Cell(17,2).Range.Text = "72×71="
$t.Cell(17,3).Range.Text = "98×39="
$t.Cell(17,4).Range.Text = "33×53="
$t.Cell(17,5).Range.Text = "48×97="
$t.Cell(18,1).Range.Text = "19×90="
$t.Cell(18,2).Range.Text = "61×98="
$t.Cell(18,3).Range.Text = "32×77="
$t.Cell(18,4).Range.Text = "69×77="
$t.Cell(18,5).Range.Text = "51×17="
$t.Cell(19,1).Range.Text = "82×83="
$t.Cell(19,2).Range.Text = "83×35="
$t.Cell(19,3).Range.Text = "86×58="
$t.Cell(19,4).Range.Text = "78×78="
$t.Cell(19,5).Range.Text = "20×43="
$t.Cell(20,1).Range.Text = "14×10="
$t.Cell(20,2).Range.Text = "58×78="
$t.Cell(20,3).Range.Text = "11×83="
$t.Cell(20,4).Range.Text = "58×69="
$t.Cell(20,5).Range.Text = "48×84="
